$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 60) with the English / Chinese / Vietnamese
# translations for "cha siu bao" (叉燒包 / bánh bao xá xíu). Fill the
# Chinese and Vietnamese cells before the English one so new shared
# strings land in the same order as the source edit.
$ws.Range("B60").Value = "叉燒包"
$ws.Range("C60").Value = "bánh bao xá xíu"
$ws.Range("A60").Value = "cha siu bao"

# Move the active selection the way it ended up after the edit in Excel.
$ws.Range("A63").Select()
